# Auto-generated Excel COM-interop script applying numeric cell updates
# sourced from the Lamia_Profits.xlsx diff (scheduled price-data refresh).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(11, 8).Value = 642.7778
$ws.Cells.Item(11, 9).Value = 642.7778
$ws.Cells.Item(11, 11).Value = 642.7778
$ws.Cells.Item(11, 13).Value = -502.7778
$ws.Cells.Item(33, 8).Value = 384.2143
$ws.Cells.Item(33, 9).Value = 299.4
$ws.Cells.Item(33, 10).Value = 596.25
$ws.Cells.Item(33, 11).Value = 299.4
$ws.Cells.Item(33, 12).Value = 596.25
$ws.Cells.Item(33, 13).Value = -70.39999999999998
$ws.Cells.Item(33, 14).Value = -1054.25
$ws.Cells.Item(38, 8).Value = 45.2
$ws.Cells.Item(38, 9).Value = 45.2
$ws.Cells.Item(38, 11).Value = 135.6
$ws.Cells.Item(38, 13).Value = 236.4
$ws.Cells.Item(39, 8).Value = 192.93103
$ws.Cells.Item(39, 9).Value = 199.67857
$ws.Cells.Item(39, 11).Value = 599.03571
$ws.Cells.Item(39, 13).Value = -303.03571
$ws.Cells.Item(45, 8).Value = 2363.6
$ws.Cells.Item(45, 10).Value = 2363.6
$ws.Cells.Item(45, 12).Value = 7090.799999999999
$ws.Cells.Item(45, 14).Value = -7474.799999999999
$ws.Cells.Item(62, 8).Value = 9572.286
$ws.Cells.Item(62, 9).Value = 8000
$ws.Cells.Item(62, 10).Value = 9834.333000000001
$ws.Cells.Item(62, 11).Value = 8000
$ws.Cells.Item(62, 12).Value = 9834.333000000001
$ws.Cells.Item(62, 13).Value = -7376
$ws.Cells.Item(62, 14).Value = -11082.333
$ws.Cells.Item(65, 8).Value = 9572.286
$ws.Cells.Item(65, 9).Value = 8000
$ws.Cells.Item(65, 10).Value = 9834.333000000001
$ws.Cells.Item(65, 11).Value = 40000
$ws.Cells.Item(65, 12).Value = 49171.665
$ws.Cells.Item(65, 13).Value = -36880
$ws.Cells.Item(65, 14).Value = -55411.665
$ws.Cells.Item(68, 8).Value = 98333.336
$ws.Cells.Item(68, 10).Value = 85000
$ws.Cells.Item(68, 12).Value = 85000
$ws.Cells.Item(68, 14).Value = -86498
$ws.Cells.Item(71, 8).Value = 98333.336
$ws.Cells.Item(71, 10).Value = 85000
$ws.Cells.Item(71, 12).Value = 255000
$ws.Cells.Item(71, 14).Value = -262488
$ws.Cells.Item(86, 8).Value = 6375.375
$ws.Cells.Item(86, 9).Value = 7501.5
$ws.Cells.Item(86, 11).Value = 7501.5
$ws.Cells.Item(86, 13).Value = -6378.5
$ws.Cells.Item(88, 8).Value = 235608.84
$ws.Cells.Item(88, 10).Value = 6196.5713
$ws.Cells.Item(88, 12).Value = 6196.5713
$ws.Cells.Item(88, 14).Value = -7008.5713
$ws.Cells.Item(89, 8).Value = 6375.375
$ws.Cells.Item(89, 9).Value = 7501.5
$ws.Cells.Item(89, 11).Value = 37507.5
$ws.Cells.Item(89, 13).Value = -31891.5
$ws.Cells.Item(91, 8).Value = 235608.84
$ws.Cells.Item(91, 10).Value = 6196.5713
$ws.Cells.Item(91, 12).Value = 6196.5713
$ws.Cells.Item(91, 13).Value = -9705.833500000001
$ws.Cells.Item(91, 14).Value = -9004.5713
$ws.Cells.Item(98, 8).Value = 250890.2
$ws.Cells.Item(98, 9).Value = 998.7
$ws.Cells.Item(98, 11).Value = 998.7
$ws.Cells.Item(98, 13).Value = 499.3
$ws.Cells.Item(106, 8).Value = 4012.2083
$ws.Cells.Item(106, 9).Value = 4064.6667
$ws.Cells.Item(106, 11).Value = 4064.6667
$ws.Cells.Item(106, 13).Value = -3433.6667
$ws.Cells.Item(113, 8).Value = 8997.4
$ws.Cells.Item(113, 9).Value = 8296.666999999999
$ws.Cells.Item(113, 10).Value = 9297.714
$ws.Cells.Item(113, 11).Value = 8296.666999999999
$ws.Cells.Item(113, 12).Value = 9297.714
$ws.Cells.Item(113, 13).Value = -5042.666999999999
$ws.Cells.Item(113, 14).Value = -15805.714
$ws.Cells.Item(116, 8).Value = 8384.177
$ws.Cells.Item(116, 9).Value = 7473.25
$ws.Cells.Item(116, 10).Value = 9193.888999999999
$ws.Cells.Item(116, 11).Value = 7473.25
$ws.Cells.Item(116, 12).Value = 9193.888999999999
$ws.Cells.Item(116, 13).Value = -4031.25
$ws.Cells.Item(116, 14).Value = -16077.889
$ws.Cells.Item(122, 8).Value = 250890.2
$ws.Cells.Item(122, 9).Value = 998.7
$ws.Cells.Item(122, 11).Value = 2996.1
$ws.Cells.Item(122, 13).Value = -546.1000000000004
$ws.Cells.Item(137, 8).Value = 40003932
$ws.Cells.Item(137, 9).Value = 62503076
$ws.Cells.Item(137, 10).Value = 5455.6665
$ws.Cells.Item(137, 11).Value = 187509228
$ws.Cells.Item(137, 12).Value = 16366.9995
$ws.Cells.Item(137, 13).Value = -187506678
$ws.Cells.Item(137, 14).Value = -21466.9995
$ws.Cells.Item(138, 8).Value = 5320.22
$ws.Cells.Item(138, 10).Value = 5259.579
$ws.Cells.Item(138, 12).Value = 15778.737
$ws.Cells.Item(138, 14).Value = -26058.737
$ws.Cells.Item(139, 8).Value = 68888
$ws.Cells.Item(139, 10).Value = 68888
$ws.Cells.Item(139, 12).Value = 68888
$ws.Cells.Item(139, 14).Value = -79168

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 7522.7646
$ws.Cells.Item(2, 9).Value = 1853.909
$ws.Cells.Item(2, 11).Value = 1853.909
$ws.Cells.Item(2, 13).Value = -1740.909
$ws.Cells.Item(5, 8).Value = 341.1111
$ws.Cells.Item(5, 9).Value = 340.83334
$ws.Cells.Item(5, 10).Value = 341.66666
$ws.Cells.Item(5, 11).Value = 340.83334
$ws.Cells.Item(5, 12).Value = 341.66666
$ws.Cells.Item(5, 13).Value = -228.83334
$ws.Cells.Item(5, 14).Value = -565.66666
$ws.Cells.Item(32, 8).Value = 9766.210999999999
$ws.Cells.Item(32, 9).Value = 7786.0757
$ws.Cells.Item(32, 10).Value = 36003
$ws.Cells.Item(32, 11).Value = 7786.0757
$ws.Cells.Item(32, 12).Value = 36003
$ws.Cells.Item(32, 13).Value = -7499.0757
$ws.Cells.Item(32, 14).Value = -36577
$ws.Cells.Item(80, 8).Value = 142632.33
$ws.Cells.Item(80, 10).Value = 141998.5
$ws.Cells.Item(80, 12).Value = 141998.5
$ws.Cells.Item(80, 14).Value = -143994.5
$ws.Cells.Item(83, 8).Value = 142632.33
$ws.Cells.Item(83, 10).Value = 141998.5
$ws.Cells.Item(83, 12).Value = 425995.5
$ws.Cells.Item(83, 14).Value = -435979.5
$ws.Cells.Item(116, 8).Value = 7522.7646
$ws.Cells.Item(116, 9).Value = 1853.909
$ws.Cells.Item(116, 11).Value = 1853.909
$ws.Cells.Item(116, 13).Value = 440.0909999999999
$ws.Cells.Item(122, 8).Value = 2931.7715
$ws.Cells.Item(122, 9).Value = 2304.96
$ws.Cells.Item(122, 11).Value = 6914.88
$ws.Cells.Item(122, 13).Value = -4464.88
$ws.Cells.Item(129, 8).Value = 60000
$ws.Cells.Item(129, 10).Value = 60000
$ws.Cells.Item(129, 12).Value = 60000
$ws.Cells.Item(129, 14).Value = -70000
$ws.Cells.Item(132, 8).Value = 3710.4285
$ws.Cells.Item(132, 9).Value = 3110.6924
$ws.Cells.Item(132, 10).Value = 11507
$ws.Cells.Item(132, 11).Value = 9332.0772
$ws.Cells.Item(132, 12).Value = 34521
$ws.Cells.Item(132, 13).Value = -6802.0772
$ws.Cells.Item(132, 14).Value = -39581

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 7522.7646
$ws.Cells.Item(3, 9).Value = 1853.909
$ws.Cells.Item(3, 11).Value = 1853.909
$ws.Cells.Item(3, 13).Value = -1739.909
$ws.Cells.Item(4, 8).Value = 341.1111
$ws.Cells.Item(4, 9).Value = 340.83334
$ws.Cells.Item(4, 10).Value = 341.66666
$ws.Cells.Item(4, 11).Value = 340.83334
$ws.Cells.Item(4, 12).Value = 341.66666
$ws.Cells.Item(4, 13).Value = -225.83334
$ws.Cells.Item(4, 14).Value = -571.66666
$ws.Cells.Item(22, 8).Value = 937.25
$ws.Cells.Item(22, 9).Value = 399.75
$ws.Cells.Item(22, 10).Value = 1474.75
$ws.Cells.Item(22, 11).Value = 399.75
$ws.Cells.Item(22, 12).Value = 1474.75
$ws.Cells.Item(22, 13).Value = -226.75
$ws.Cells.Item(22, 14).Value = -1820.75
$ws.Cells.Item(86, 8).Value = 4523.5
$ws.Cells.Item(86, 9).Value = 3706.5625
$ws.Cells.Item(86, 10).Value = 5830.6
$ws.Cells.Item(86, 11).Value = 3706.5625
$ws.Cells.Item(86, 12).Value = 5830.6
$ws.Cells.Item(86, 13).Value = -2583.5625
$ws.Cells.Item(86, 14).Value = -8076.6
$ws.Cells.Item(89, 8).Value = 4523.5
$ws.Cells.Item(89, 9).Value = 3706.5625
$ws.Cells.Item(89, 10).Value = 5830.6
$ws.Cells.Item(89, 11).Value = 18532.8125
$ws.Cells.Item(89, 12).Value = 29153
$ws.Cells.Item(89, 13).Value = -12916.8125
$ws.Cells.Item(89, 14).Value = -40385
$ws.Cells.Item(94, 8).Value = 2548.8845
$ws.Cells.Item(94, 9).Value = 2471.6365
$ws.Cells.Item(94, 11).Value = 2471.6365
$ws.Cells.Item(94, 13).Value = -2020.6365
$ws.Cells.Item(99, 8).Value = 3791.818
$ws.Cells.Item(99, 9).Value = 2066.5
$ws.Cells.Item(99, 10).Value = 5862.2
$ws.Cells.Item(99, 11).Value = 2066.5
$ws.Cells.Item(99, 12).Value = 5862.2
$ws.Cells.Item(99, 13).Value = -568.5
$ws.Cells.Item(99, 14).Value = -8858.200000000001
$ws.Cells.Item(105, 8).Value = 18396.053
$ws.Cells.Item(105, 9).Value = 19370.363
$ws.Cells.Item(105, 11).Value = 19370.363
$ws.Cells.Item(105, 13).Value = -17623.363
$ws.Cells.Item(110, 8).Value = 66473.336
$ws.Cells.Item(110, 10).Value = 66473.336
$ws.Cells.Item(110, 12).Value = 66473.336
$ws.Cells.Item(110, 14).Value = -74653.336
$ws.Cells.Item(138, 8).Value = 69992.8
$ws.Cells.Item(138, 10).Value = 69992.8
$ws.Cells.Item(138, 12).Value = 69992.8
$ws.Cells.Item(138, 14).Value = -80272.8

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(16, 8).Value = 2787.3333
$ws.Cells.Item(16, 9).Value = 1712.1666
$ws.Cells.Item(16, 10).Value = 3862.5
$ws.Cells.Item(16, 11).Value = 1712.1666
$ws.Cells.Item(16, 12).Value = 3862.5
$ws.Cells.Item(16, 13).Value = -1425.1666
$ws.Cells.Item(16, 14).Value = -4436.5
$ws.Cells.Item(17, 8).Value = 500
$ws.Cells.Item(17, 9).Value = 500
$ws.Cells.Item(17, 11).Value = 500
$ws.Cells.Item(17, 13).Value = -326
$ws.Cells.Item(22, 9).Value = 3255.3333
$ws.Cells.Item(22, 10).Value = 10257.5
$ws.Cells.Item(22, 11).Value = 3255.3333
$ws.Cells.Item(22, 12).Value = 10257.5
$ws.Cells.Item(22, 13).Value = -2905.3333
$ws.Cells.Item(22, 14).Value = -10957.5
$ws.Cells.Item(31, 8).Value = 43090.793
$ws.Cells.Item(31, 9).Value = 11991.462
$ws.Cells.Item(31, 10).Value = 68359
$ws.Cells.Item(31, 11).Value = 11991.462
$ws.Cells.Item(31, 12).Value = 68359
$ws.Cells.Item(31, 13).Value = -11696.462
$ws.Cells.Item(31, 14).Value = -68949
$ws.Cells.Item(34, 8).Value = 43090.793
$ws.Cells.Item(34, 9).Value = 11991.462
$ws.Cells.Item(34, 10).Value = 68359
$ws.Cells.Item(34, 11).Value = 11991.462
$ws.Cells.Item(34, 12).Value = 68359
$ws.Cells.Item(34, 13).Value = -11789.462
$ws.Cells.Item(34, 14).Value = -68763
$ws.Cells.Item(41, 8).Value = 16686.334
$ws.Cells.Item(50, 8).Value = 26249.25
$ws.Cells.Item(50, 10).Value = 29999
$ws.Cells.Item(50, 12).Value = 29999
$ws.Cells.Item(50, 14).Value = -31249
$ws.Cells.Item(51, 8).Value = 39999.723
$ws.Cells.Item(51, 10).Value = 39999.766
$ws.Cells.Item(51, 12).Value = 39999.766
$ws.Cells.Item(51, 14).Value = -41471.766
$ws.Cells.Item(58, 8).Value = 4823.905
$ws.Cells.Item(58, 9).Value = 1800.9
$ws.Cells.Item(58, 10).Value = 7572.091
$ws.Cells.Item(58, 11).Value = 1800.9
$ws.Cells.Item(58, 12).Value = 7572.091
$ws.Cells.Item(58, 13).Value = -1597.9
$ws.Cells.Item(58, 14).Value = -7978.091
$ws.Cells.Item(59, 8).Value = 28109.666
$ws.Cells.Item(59, 9).Value = 22165.834
$ws.Cells.Item(59, 10).Value = 39997.332
$ws.Cells.Item(59, 11).Value = 22165.834
$ws.Cells.Item(59, 12).Value = 39997.332
$ws.Cells.Item(59, 13).Value = -21020.834
$ws.Cells.Item(59, 14).Value = -42287.332
$ws.Cells.Item(60, 8).Value = 35021.043
$ws.Cells.Item(60, 10).Value = 39999.2
$ws.Cells.Item(60, 12).Value = 39999.2
$ws.Cells.Item(60, 14).Value = -41021.2
$ws.Cells.Item(61, 8).Value = 39999.723
$ws.Cells.Item(61, 10).Value = 39999.766
$ws.Cells.Item(61, 12).Value = 39999.766
$ws.Cells.Item(61, 14).Value = -40695.766
$ws.Cells.Item(68, 8).Value = 48570
$ws.Cells.Item(68, 10).Value = 48570
$ws.Cells.Item(68, 12).Value = 48570
$ws.Cells.Item(68, 14).Value = -50068
$ws.Cells.Item(71, 8).Value = 48570
$ws.Cells.Item(71, 10).Value = 48570
$ws.Cells.Item(71, 12).Value = 145710
$ws.Cells.Item(71, 14).Value = -153198
$ws.Cells.Item(74, 8).Value = 66525.336
$ws.Cells.Item(77, 8).Value = 66525.336
$ws.Cells.Item(82, 8).Value = 34980.5
$ws.Cells.Item(82, 10).Value = 34980.5
$ws.Cells.Item(82, 12).Value = 34980.5
$ws.Cells.Item(82, 14).Value = -35702.5
$ws.Cells.Item(85, 8).Value = 34980.5
$ws.Cells.Item(85, 10).Value = 34980.5
$ws.Cells.Item(85, 12).Value = 34980.5
$ws.Cells.Item(85, 14).Value = -37476.5
$ws.Cells.Item(105, 8).Value = 1838.1428
$ws.Cells.Item(105, 9).Value = 663.4167
$ws.Cells.Item(105, 11).Value = 663.4167
$ws.Cells.Item(105, 13).Value = 1083.5833
$ws.Cells.Item(107, 8).Value = 2493.76
$ws.Cells.Item(107, 9).Value = 1969.2778
$ws.Cells.Item(107, 10).Value = 3842.4285
$ws.Cells.Item(107, 11).Value = 1969.2778
$ws.Cells.Item(107, 12).Value = 3842.4285
$ws.Cells.Item(107, 13).Value = -49.27780000000007
$ws.Cells.Item(107, 14).Value = -7682.4285
$ws.Cells.Item(113, 8).Value = 2787.3333
$ws.Cells.Item(113, 9).Value = 1712.1666
$ws.Cells.Item(113, 10).Value = 3862.5
$ws.Cells.Item(113, 11).Value = 1712.1666
$ws.Cells.Item(113, 12).Value = 3862.5
$ws.Cells.Item(113, 13).Value = 457.8334
$ws.Cells.Item(113, 14).Value = -8202.5
$ws.Cells.Item(122, 8).Value = 5602.9443
$ws.Cells.Item(122, 9).Value = 1955.8889
$ws.Cells.Item(122, 10).Value = 9250
$ws.Cells.Item(122, 11).Value = 5867.6667
$ws.Cells.Item(122, 12).Value = 27750
$ws.Cells.Item(122, 13).Value = -3417.6667
$ws.Cells.Item(122, 14).Value = -32650
$ws.Cells.Item(132, 8).Value = 3880.74
$ws.Cells.Item(132, 9).Value = 3407.4048
$ws.Cells.Item(132, 10).Value = 6365.75
$ws.Cells.Item(132, 11).Value = 10222.2144
$ws.Cells.Item(132, 12).Value = 19097.25
$ws.Cells.Item(132, 13).Value = -7692.214399999999
$ws.Cells.Item(132, 14).Value = -24157.25
$ws.Cells.Item(134, 8).Value = 2245.1482
$ws.Cells.Item(134, 9).Value = 1742.5834
$ws.Cells.Item(134, 11).Value = 5227.7502
$ws.Cells.Item(134, 13).Value = -2692.7502
$ws.Cells.Item(136, 8).Value = 4823.905
$ws.Cells.Item(136, 9).Value = 1800.9
$ws.Cells.Item(136, 10).Value = 7572.091
$ws.Cells.Item(136, 11).Value = 5402.700000000001
$ws.Cells.Item(136, 12).Value = 22716.273
$ws.Cells.Item(136, 13).Value = -2852.700000000001
$ws.Cells.Item(136, 14).Value = -27816.273

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(5, 8).Value = 2382291.8
$ws.Cells.Item(5, 9).Value = 1345.1
$ws.Cells.Item(5, 10).Value = 8334658
$ws.Cells.Item(5, 11).Value = 4035.3
$ws.Cells.Item(5, 12).Value = 25003974
$ws.Cells.Item(5, 13).Value = -3923.3
$ws.Cells.Item(5, 14).Value = -25004198
$ws.Cells.Item(7, 8).Value = 116841.14
$ws.Cells.Item(7, 9).Value = 79372.5
$ws.Cells.Item(7, 10).Value = 166799.33
$ws.Cells.Item(7, 11).Value = 238117.5
$ws.Cells.Item(7, 12).Value = 500397.99
$ws.Cells.Item(7, 13).Value = -238005.5
$ws.Cells.Item(7, 14).Value = -500621.99
$ws.Cells.Item(39, 8).Value = 2999.75
$ws.Cells.Item(39, 10).Value = 2666.3333
$ws.Cells.Item(39, 12).Value = 7998.999899999999
$ws.Cells.Item(39, 14).Value = -8586.999899999999
$ws.Cells.Item(44, 8).Value = 142965.14
$ws.Cells.Item(44, 10).Value = 166792.17
$ws.Cells.Item(44, 12).Value = 500376.51
$ws.Cells.Item(44, 14).Value = -501172.51
$ws.Cells.Item(107, 8).Value = 3472642.5
$ws.Cells.Item(107, 9).Value = 472.75
$ws.Cells.Item(107, 11).Value = 1418.25
$ws.Cells.Item(107, 13).Value = 501.75
$ws.Cells.Item(113, 8).Value = 2045
$ws.Cells.Item(113, 10).Value = 2045
$ws.Cells.Item(113, 12).Value = 6135
$ws.Cells.Item(113, 14).Value = -10475
$ws.Cells.Item(127, 8).Value = 936
$ws.Cells.Item(127, 10).Value = 936
$ws.Cells.Item(127, 12).Value = 2808
$ws.Cells.Item(127, 14).Value = -12728
$ws.Cells.Item(135, 8).Value = 2382291.8
$ws.Cells.Item(135, 9).Value = 1345.1
$ws.Cells.Item(135, 10).Value = 8334658
$ws.Cells.Item(135, 11).Value = 12105.9
$ws.Cells.Item(135, 12).Value = 75011922
$ws.Cells.Item(135, 13).Value = -9570.9
$ws.Cells.Item(135, 14).Value = -75016992
$ws.Cells.Item(138, 8).Value = 2395
$ws.Cells.Item(138, 9).Value = 1472
$ws.Cells.Item(138, 10).Value = 3933.3333
$ws.Cells.Item(138, 11).Value = 4416
$ws.Cells.Item(138, 12).Value = 11799.9999
$ws.Cells.Item(138, 13).Value = 724
$ws.Cells.Item(138, 14).Value = -22079.9999

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 8).Value = 439.44446
$ws.Cells.Item(2, 10).Value = 1015
$ws.Cells.Item(2, 12).Value = 1015
$ws.Cells.Item(2, 14).Value = -1241
$ws.Cells.Item(102, 8).Value = 2313.8965
$ws.Cells.Item(102, 9).Value = 944.5
$ws.Cells.Item(102, 10).Value = 5357
$ws.Cells.Item(102, 11).Value = 944.5
$ws.Cells.Item(102, 12).Value = 5357
$ws.Cells.Item(102, 13).Value = 677.5
$ws.Cells.Item(102, 14).Value = -8601
$ws.Cells.Item(113, 8).Value = 3328.087
$ws.Cells.Item(113, 9).Value = 2446.7273
$ws.Cells.Item(113, 11).Value = 2446.7273
$ws.Cells.Item(113, 13).Value = -276.7273
$ws.Cells.Item(126, 8).Value = 4219.1777
$ws.Cells.Item(126, 9).Value = 4376.7393
$ws.Cells.Item(126, 11).Value = 13130.2179
$ws.Cells.Item(126, 13).Value = -10660.2179
$ws.Cells.Item(132, 8).Value = 2846.8096
$ws.Cells.Item(132, 9).Value = 1960.8125
$ws.Cells.Item(132, 10).Value = 5682
$ws.Cells.Item(132, 11).Value = 5882.4375
$ws.Cells.Item(132, 12).Value = 17046
$ws.Cells.Item(132, 13).Value = -3352.4375
$ws.Cells.Item(132, 14).Value = -22106

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 8240.704
$ws.Cells.Item(7, 9).Value = 6112.45
$ws.Cells.Item(7, 11).Value = 6112.45
$ws.Cells.Item(7, 13).Value = -6000.45
$ws.Cells.Item(22, 8).Value = 166675070
$ws.Cells.Item(22, 10).Value = 16000
$ws.Cells.Item(22, 12).Value = 16000
$ws.Cells.Item(22, 14).Value = -16590
$ws.Cells.Item(27, 8).Value = 166675070
$ws.Cells.Item(27, 10).Value = 16000
$ws.Cells.Item(27, 12).Value = 16000
$ws.Cells.Item(27, 14).Value = -16214
$ws.Cells.Item(46, 8).Value = 5424.5713
$ws.Cells.Item(46, 9).Value = 2485
$ws.Cells.Item(46, 10).Value = 6600.4
$ws.Cells.Item(46, 11).Value = 2485
$ws.Cells.Item(46, 12).Value = 6600.4
$ws.Cells.Item(46, 13).Value = -2297
$ws.Cells.Item(46, 14).Value = -6976.4
$ws.Cells.Item(55, 8).Value = 5885969.5
$ws.Cells.Item(55, 9).Value = 16666916
$ws.Cells.Item(55, 10).Value = 5453
$ws.Cells.Item(55, 11).Value = 16666916
$ws.Cells.Item(55, 12).Value = 5453
$ws.Cells.Item(55, 13).Value = -16666743
$ws.Cells.Item(55, 14).Value = -5799
$ws.Cells.Item(100, 8).Value = 5659.421
$ws.Cells.Item(100, 9).Value = 4709.8667
$ws.Cells.Item(100, 11).Value = 4709.8667
$ws.Cells.Item(100, 13).Value = -4168.8667
$ws.Cells.Item(121, 8).Value = 47525.332
$ws.Cells.Item(121, 10).Value = 47525.332
$ws.Cells.Item(121, 12).Value = 47525.332
$ws.Cells.Item(121, 14).Value = -51019.332
$ws.Cells.Item(122, 8).Value = 243486.36
$ws.Cells.Item(122, 9).Value = 450820.88
$ws.Cells.Item(122, 11).Value = 1352462.64
$ws.Cells.Item(122, 13).Value = -1350012.64
$ws.Cells.Item(126, 8).Value = 8240.704
$ws.Cells.Item(126, 9).Value = 6112.45
$ws.Cells.Item(126, 11).Value = 18337.35
$ws.Cells.Item(126, 13).Value = -15867.35
$ws.Cells.Item(132, 8).Value = 4850.1
$ws.Cells.Item(132, 9).Value = 4425.25
$ws.Cells.Item(132, 10).Value = 5487.375
$ws.Cells.Item(132, 11).Value = 13275.75
$ws.Cells.Item(132, 12).Value = 16462.125
$ws.Cells.Item(132, 13).Value = -10745.75
$ws.Cells.Item(132, 14).Value = -21522.125
$ws.Cells.Item(140, 8).Value = 99000
$ws.Cells.Item(140, 10).Value = 99000
$ws.Cells.Item(140, 12).Value = 99000
$ws.Cells.Item(140, 14).Value = -109360

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(62, 8).Value = 7670.2
$ws.Cells.Item(62, 9).Value = 4700
$ws.Cells.Item(62, 10).Value = 8000.222
$ws.Cells.Item(62, 11).Value = 4700
$ws.Cells.Item(62, 12).Value = 8000.222
$ws.Cells.Item(62, 13).Value = -4076
$ws.Cells.Item(62, 14).Value = -9248.222
$ws.Cells.Item(65, 8).Value = 7670.2
$ws.Cells.Item(65, 9).Value = 4700
$ws.Cells.Item(65, 10).Value = 8000.222
$ws.Cells.Item(65, 11).Value = 23500
$ws.Cells.Item(65, 12).Value = 40001.11
$ws.Cells.Item(65, 13).Value = -20380
$ws.Cells.Item(65, 14).Value = -46241.11
$ws.Cells.Item(92, 8).Value = 25000
$ws.Cells.Item(92, 10).Value = 25000
$ws.Cells.Item(92, 12).Value = 25000
$ws.Cells.Item(92, 14).Value = -29992
$ws.Cells.Item(122, 8).Value = 2596.182
$ws.Cells.Item(122, 9).Value = 1532.25
$ws.Cells.Item(122, 10).Value = 5433.3335
$ws.Cells.Item(122, 11).Value = 4596.75
$ws.Cells.Item(122, 12).Value = 16300.0005
$ws.Cells.Item(122, 13).Value = -2146.75
$ws.Cells.Item(122, 14).Value = -21200.0005
$ws.Cells.Item(126, 8).Value = 3174.5151
$ws.Cells.Item(126, 9).Value = 1876.3077
$ws.Cells.Item(126, 10).Value = 7996.4287
$ws.Cells.Item(126, 11).Value = 5628.9231
$ws.Cells.Item(126, 12).Value = 23989.2861
$ws.Cells.Item(126, 13).Value = -3158.9231
$ws.Cells.Item(126, 14).Value = -28929.2861
$ws.Cells.Item(132, 8).Value = 4174.2354
$ws.Cells.Item(132, 9).Value = 2304.3845
$ws.Cells.Item(132, 10).Value = 10251.25
$ws.Cells.Item(132, 11).Value = 6913.1535
$ws.Cells.Item(132, 12).Value = 30753.75
$ws.Cells.Item(132, 13).Value = -4383.1535
$ws.Cells.Item(132, 14).Value = -35813.75
$ws.Cells.Item(136, 8).Value = 3365.5
$ws.Cells.Item(136, 9).Value = 2445.5881
$ws.Cells.Item(136, 10).Value = 19004
$ws.Cells.Item(136, 11).Value = 7336.7643
$ws.Cells.Item(136, 12).Value = 57012
$ws.Cells.Item(136, 13).Value = -4786.7643
$ws.Cells.Item(136, 14).Value = -62112
